$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H54").Value = 5000
$ws.Range("J54").Value = 5000
$ws.Range("L54").Value = 5000
$ws.Range("N54").Value = -5972

$ws.Range("H95").Value = 8522.223
$ws.Range("J95").Value = 8522.223
$ws.Range("L95").Value = 8522.223
$ws.Range("N95").Value = -14014.223

$ws.Range("H97").Value = 4925
$ws.Range("J97").Value = 4925
$ws.Range("L97").Value = 14775
$ws.Range("N97").Value = -15767

$ws.Range("H113").Value = 50002624
$ws.Range("J113").Value = 3499.5
$ws.Range("L113").Value = 3499.5
$ws.Range("N113").Value = -10007.5

$ws.Range("H121").Value = 768.84
$ws.Range("J121").Value = 791.2917
$ws.Range("L121").Value = 2373.8751
$ws.Range("N121").Value = -5867.8751

$ws.Range("H132").Value = 6674424
$ws.Range("I132").Value = 7581118.5
$ws.Range("K132").Value = 22743355.5
$ws.Range("M132").Value = -22740825.5

$ws.Range("H138").Value = 3708.6455
$ws.Range("I138").Value = 2332.6924
$ws.Range("J138").Value = 3979.6667
$ws.Range("K138").Value = 6998.0772
$ws.Range("L138").Value = 11939.0001
$ws.Range("M138").Value = -1858.0772
$ws.Range("N138").Value = -22219.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24798.445
$ws.Range("I32").Value = 22273.312
$ws.Range("J32").Value = 44999.5
$ws.Range("K32").Value = 22273.312
$ws.Range("L32").Value = 44999.5
$ws.Range("M32").Value = -21986.312
$ws.Range("N32").Value = -45573.5

$ws.Range("H61").Value = 45456188
$ws.Range("I61").Value = 55556964
$ws.Range("J61").Value = 2710.5
$ws.Range("K61").Value = 55556964
$ws.Range("L61").Value = 2710.5
$ws.Range("M61").Value = -55556752
$ws.Range("N61").Value = -3134.5

$ws.Range("H136").Value = 45456188
$ws.Range("I136").Value = 55556964
$ws.Range("J136").Value = 2710.5
$ws.Range("K136").Value = 166670892
$ws.Range("L136").Value = 8131.5
$ws.Range("M136").Value = -166668342
$ws.Range("N136").Value = -13231.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 20751
$ws.Range("I82").Value = 12564.25
$ws.Range("J82").Value = 31666.666
$ws.Range("K82").Value = 12564.25
$ws.Range("L82").Value = 31666.666
$ws.Range("M82").Value = -12181.25
$ws.Range("N82").Value = -32432.666

$ws.Range("H85").Value = 20751
$ws.Range("I85").Value = 12564.25
$ws.Range("J85").Value = 31666.666
$ws.Range("K85").Value = 12564.25
$ws.Range("L85").Value = 31666.666
$ws.Range("M85").Value = -11238.25
$ws.Range("N85").Value = -34318.666

$ws.Range("H94").Value = 10000450
$ws.Range("I94").Value = 10417115
$ws.Range("J94").Value = 479
$ws.Range("K94").Value = 10417115
$ws.Range("L94").Value = 479
$ws.Range("M94").Value = -10416664
$ws.Range("N94").Value = -1381

$ws.Range("H138").Value = 48766.668
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 48766.668
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 48766.668
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -59046.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1638.0597
$ws.Range("I31").Value = 1361.3793
$ws.Range("J31").Value = 3421.111
$ws.Range("K31").Value = 1361.3793
$ws.Range("L31").Value = 3421.111
$ws.Range("M31").Value = -1066.3793
$ws.Range("N31").Value = -4011.111

$ws.Range("H34").Value = 1638.0597
$ws.Range("I34").Value = 1361.3793
$ws.Range("J34").Value = 3421.111
$ws.Range("K34").Value = 1361.3793
$ws.Range("L34").Value = 3421.111
$ws.Range("M34").Value = -1159.3793
$ws.Range("N34").Value = -3825.111

$ws.Range("H141").Value = 534023
$ws.Range("J141").Value = 534023
$ws.Range("L141").Value = 534023
$ws.Range("N141").Value = -544383

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 13651.615
$ws.Range("I3").Value = 9357.143
$ws.Range("K3").Value = 28071.429
$ws.Range("M3").Value = -27959.429

$ws.Range("H107").Value = 8020.3125
$ws.Range("I107").Value = 391.66666
$ws.Range("J107").Value = 12597.5
$ws.Range("K107").Value = 1174.99998
$ws.Range("L107").Value = 37792.5
$ws.Range("M107").Value = 745.0000199999999
$ws.Range("N107").Value = -41632.5

$ws.Range("H122").Value = 1236.6666
$ws.Range("I122").Value = 872.6667
$ws.Range("J122").Value = 1600.6666
$ws.Range("K122").Value = 7854.0003
$ws.Range("L122").Value = 14405.9994
$ws.Range("M122").Value = -5404.0003
$ws.Range("N122").Value = -19305.9994

$ws.Range("H131").Value = 23846324
$ws.Range("J131").Value = 38615.05
$ws.Range("L131").Value = 115845.15
$ws.Range("N131").Value = -125925.15

$ws.Range("H140").Value = 29982.621
$ws.Range("I140").Value = 80004.46000000001
$ws.Range("J140").Value = 2887.4583
$ws.Range("K140").Value = 240013.38
$ws.Range("L140").Value = 8662.374899999999
$ws.Range("M140").Value = -234833.38
$ws.Range("N140").Value = -19022.3749

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1335.091
$ws.Range("I113").Value = 1374.7778
$ws.Range("J113").Value = 1156.5
$ws.Range("K113").Value = 1374.7778
$ws.Range("L113").Value = 1156.5
$ws.Range("M113").Value = 795.2221999999999
$ws.Range("N113").Value = -5496.5

$ws.Range("H122").Value = 1557.9375
$ws.Range("I122").Value = 1494.4615
$ws.Range("J122").Value = 1833
$ws.Range("K122").Value = 4483.3845
$ws.Range("L122").Value = 5499
$ws.Range("M122").Value = -2033.3845
$ws.Range("N122").Value = -10399

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 995
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()

$ws.Range("H22").Value = 1700.1666
$ws.Range("J22").Value = 1176
$ws.Range("L22").Value = 1176
$ws.Range("N22").Value = -1766

$ws.Range("H27").Value = 1700.1666
$ws.Range("J27").Value = 1176
$ws.Range("L27").Value = 1176
$ws.Range("N27").Value = -1390

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4024.487
$ws.Range("I132").Value = 4107.643
$ws.Range("J132").Value = 3812.818
$ws.Range("K132").Value = 12322.929
$ws.Range("L132").Value = 11438.454
$ws.Range("M132").Value = -9792.929
$ws.Range("N132").Value = -16498.454
